# Fix element ordering inside <w:rPr> for several character styles in
# styles.xml so that <w:b/> / <w:i/> precede <w:color/>, matching the
# sequence required by wml.xsd (CT_RPr schema order).
#
# Re-assigning Font.Bold / Font.Italic to their current value forces the
# style's run properties to be re-serialized in the correct schema order
# without changing any visible formatting.

$d = $word.ActiveDocument

$boldOnly = @("KeywordTok", "ImportTok", "ControlFlowTok", "AlertTok", "ErrorTok")
foreach ($styleName in $boldOnly) {
    $s = $d.Styles.Item($styleName)
    $s.Font.Bold = $true
}

$italicOnly = @("CommentTok", "DocumentationTok")
foreach ($styleName in $italicOnly) {
    $s = $d.Styles.Item($styleName)
    $s.Font.Italic = $true
}

$boldAndItalic = @("AnnotationTok", "CommentVarTok", "InformationTok", "WarningTok")
foreach ($styleName in $boldAndItalic) {
    $s = $d.Styles.Item($styleName)
    $s.Font.Bold = $true
    $s.Font.Italic = $true
}
